$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("A25").Value = 112306136
$ws.Range("B25").Value = 90818
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 4368
$ws.Range("F25").Value = "Dofttaggsvamp"
$ws.Range("G25").Value = "Hydnellum suaveolens"
$ws.Range("H25").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("I25").NumberFormat = "General"
$ws.Range("J25").Value = "fruktkroppar"
$ws.Range("K25").NumberFormat = "General"
$ws.Range("N25").NumberFormat = "General"
$ws.Range("P25").Value = "Bässe söder om stigen, Hls"
$ws.Range("Q25").Value = 599416
$ws.Range("R25").Value = 6820643
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = "Gävleborg"
$ws.Range("U25").Value = "Hudiksvall"
$ws.Range("V25").Value = "Hälsingland"
$ws.Range("W25").Value = "Enånger"
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("Y25").Value = "2023-09-18"
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("AA25").Value = "2023-09-18"
$ws.Range("AC25").Value = "Förekommer på flera platser"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AF25").NumberFormat = "General"
$ws.Range("AG25").Value = $false
$ws.Range("AH25").Value = "Skogsmark"
$ws.Range("AI25").Value = "Barrblandskog kalkpåverkad."
$ws.Range("AT25").NumberFormat = "General"
$ws.Range("AW25").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AX25").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AY25").NumberFormat = "General"

# Row 26
$ws.Range("A26").Value = 112306179
$ws.Range("B26").Value = 89100
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 5754
$ws.Range("F26").Value = "Gultoppig fingersvamp"
$ws.Range("G26").Value = "Ramaria testaceoflava"
$ws.Range("H26").Value = "(Bres.) Corner"
$ws.Range("I26").NumberFormat = "General"
$ws.Range("J26").Value = "fruktkroppar"
$ws.Range("K26").NumberFormat = "General"
$ws.Range("N26").NumberFormat = "General"
$ws.Range("P26").Value = "Storåsens sydsluttning  söder om myren, Hls"
$ws.Range("Q26").Value = 599447
$ws.Range("R26").Value = 6820628
$ws.Range("S26").Value = 10
$ws.Range("T26").Value = "Gävleborg"
$ws.Range("U26").Value = "Hudiksvall"
$ws.Range("V26").Value = "Hälsingland"
$ws.Range("W26").Value = "Enånger"
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("Y26").Value = "2023-09-18"
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("AA26").Value = "2023-09-18"
$ws.Range("AC26").Value = "Där stigen delar sig ned mot myren"
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AF26").NumberFormat = "General"
$ws.Range("AG26").Value = $false
$ws.Range("AH26").Value = "Skogsmark"
$ws.Range("AI26").Value = "Barrblandskog kalkpåverkad."
$ws.Range("AT26").NumberFormat = "General"
$ws.Range("AW26").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AX26").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AY26").NumberFormat = "General"

# Row 27
$ws.Range("A27").Value = 112306159
$ws.Range("B27").Value = 90796
$ws.Range("C27").Value = "Ovaliderad"
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 4363
$ws.Range("F27").Value = "Zontaggsvamp"
$ws.Range("G27").Value = "Hydnellum concrescens"
$ws.Range("H27").Value = "(Pers.) Banker"
$ws.Range("I27").NumberFormat = "General"
$ws.Range("J27").NumberFormat = "General"
$ws.Range("K27").NumberFormat = "General"
$ws.Range("N27").NumberFormat = "General"
$ws.Range("P27").Value = "Storåsens sydsluttning  söder om myren, Hls"
$ws.Range("Q27").Value = 599447
$ws.Range("R27").Value = 6820628
$ws.Range("S27").Value = 10
$ws.Range("T27").Value = "Gävleborg"
$ws.Range("U27").Value = "Hudiksvall"
$ws.Range("V27").Value = "Hälsingland"
$ws.Range("W27").Value = "Enånger"
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("Y27").Value = "2023-09-18"
$ws.Range("AA27").NumberFormat = "@"
$ws.Range("AA27").Value = "2023-09-18"
$ws.Range("AC27").Value = "Efter stigen i början nära vändplatsen"
$ws.Range("AD27").Value = $false
$ws.Range("AE27").Value = $false
$ws.Range("AF27").NumberFormat = "General"
$ws.Range("AG27").Value = $false
$ws.Range("AT27").NumberFormat = "General"
$ws.Range("AW27").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AX27").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AY27").NumberFormat = "General"

# Row 28
$ws.Range("A28").Value = 112306119
$ws.Range("B28").Value = 90812
$ws.Range("C28").Value = "Ovaliderad"
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 4366
$ws.Range("F28").Value = "Skarp dropptaggsvamp"
$ws.Range("G28").Value = "Hydnellum peckii"
$ws.Range("H28").Value = "Banker"
$ws.Range("I28").NumberFormat = "General"
$ws.Range("J28").Value = "fruktkroppar"
$ws.Range("K28").NumberFormat = "General"
$ws.Range("N28").NumberFormat = "General"
$ws.Range("P28").Value = "Bässe söder om stigen, Hls"
$ws.Range("Q28").Value = 599416
$ws.Range("R28").Value = 6820643
$ws.Range("S28").Value = 10
$ws.Range("T28").Value = "Gävleborg"
$ws.Range("U28").Value = "Hudiksvall"
$ws.Range("V28").Value = "Hälsingland"
$ws.Range("W28").Value = "Enånger"
$ws.Range("Y28").NumberFormat = "@"
$ws.Range("Y28").Value = "2023-09-18"
$ws.Range("AA28").NumberFormat = "@"
$ws.Range("AA28").Value = "2023-09-18"
$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AF28").NumberFormat = "General"
$ws.Range("AG28").Value = $false
$ws.Range("AT28").NumberFormat = "General"
$ws.Range("AW28").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AX28").Value = "Alf Bjarne Roland Pallin"
$ws.Range("AY28").NumberFormat = "General"

Write-Host "done"